$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-13 hold mushroom observation records that were cyclically
# re-shuffled in the source export: the record that used to live on row 10
# is now on row 9, row 11's record moved to row 10, row 12's record moved
# to row 11, row 13's record moved to row 12, and the record that used to
# be on row 9 moved down to row 13.
#
# Reading values back out of the sheet is not reliable in this runtime, so
# the literal field values (taken from the original rows) are written
# directly to the cells that actually change.
#
# Helper: write a text value to a cell while preventing Excel from
# auto-converting date-like or numeric-like strings (e.g. "2023-08-12",
# "3") into a date serial number / real number, and without leaving any
# leftover cell formatting behind.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 9 (now holds what used to be row 10's record)
$ws.Range("A9").Value = 111611146
$ws.Range("B9").Value = 88630
Set-TextValue $ws.Range("D9") "LC"
$ws.Range("E9").Value = 4823
Set-TextValue $ws.Range("F9") "Hasselsopp"
Set-TextValue $ws.Range("G9") "Leccinellum pseudoscabrum"
Set-TextValue $ws.Range("H9") "(Kallenb.) Mikšík"

# Row 10 (now holds what used to be row 11's record)
$ws.Range("A10").Value = 111611158
$ws.Range("B10").Value = 86021
$ws.Range("E10").Value = 4037
Set-TextValue $ws.Range("F10") "Bolmörtsskivling"
Set-TextValue $ws.Range("G10") "Entoloma sinuatum"
Set-TextValue $ws.Range("H10") "(Bull.) P.Kumm."
$ws.Range("Q10").Value = 663128.0992466732
$ws.Range("R10").Value = 6634761.25188593
Set-TextValue $ws.Range("Y10") "2023-08-12"
Set-TextValue $ws.Range("AA10") "2023-08-12"
Set-TextValue $ws.Range("AC10") "1 ex. i lövförna under ek och hassel."
Set-TextValue $ws.Range("AX10") "Gillis Aronsson, Cajsa Björkén"

# Row 11 (now holds what used to be row 12's record)
$ws.Range("A11").Value = 111611145
$ws.Range("B11").Value = 88630
$ws.Range("E11").Value = 4823
Set-TextValue $ws.Range("F11") "Hasselsopp"
Set-TextValue $ws.Range("G11") "Leccinellum pseudoscabrum"
Set-TextValue $ws.Range("H11") "(Kallenb.) Mikšík"
$ws.Range("Q11").Value = 663143.8264147732
$ws.Range("R11").Value = 6634793.669287071
Set-TextValue $ws.Range("AC11") "1 ex. i lövförna under hassel."

# Row 12 (now holds what used to be row 13's record)
$ws.Range("A12").Value = 111611138
$ws.Range("B12").Value = 81796
$ws.Range("E12").Value = 5406
Set-TextValue $ws.Range("F12") "Gulmjölkig storskål"
Set-TextValue $ws.Range("G12") "Peziza succosa"
Set-TextValue $ws.Range("H12") "Berk."
Set-TextValue $ws.Range("I12") "3"
$ws.Range("Q12").Value = 663213.3366271106
$ws.Range("R12").Value = 6634830.464506784
Set-TextValue $ws.Range("AC12") "3 ex. på bar jord och i lövförna."

# Row 13 (now holds what used to be row 9's original record)
$ws.Range("A13").Value = 111611165
$ws.Range("B13").Value = 84741
Set-TextValue $ws.Range("D13") "NT"
$ws.Range("E13").Value = 37
Set-TextValue $ws.Range("F13") "Jättekamskivling"
Set-TextValue $ws.Range("G13") "Amanita ceciliae"
Set-TextValue $ws.Range("H13") "(Berk. & Broome) Bas"
Set-TextValue $ws.Range("I13") "1"
$ws.Range("Q13").Value = 663088.0668624006
$ws.Range("R13").Value = 6634684.960451891
Set-TextValue $ws.Range("Y13") "2023-08-11"
Set-TextValue $ws.Range("AA13") "2023-08-11"
Set-TextValue $ws.Range("AC13") "1 ex. under ek och hassel."
Set-TextValue $ws.Range("AX13") "Gillis Aronsson"
